{"js": "// Update the worksheet table of three-digit / one-digit division problems.\n// Cells are addressed positionally (row, column) rather than via a global\n// text search/replace, because some of the new values happen to equal the\n// old value of a different cell (e.g. the new value of row 0-index 8/col 0\n// is \"885\u00f74=221, 1\", which was the old value of row 0-index 12/col 2) --\n// replacing by position avoids accidentally re-matching an already-edited\n// cell.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of 0-based row index -> new cell values (columns 0..4), matching the\n// rows in the table that actually contain the \"a\u00f7b=c, d\" answers.\nconst updates = {\n  0: [\"382\u00f72=191, 0\", \"529\u00f78=66, 1\", \"704\u00f72=352, 0\", \"424\u00f78=53, 0\", \"145\u00f76=24, 1\"],\n  4: [\"295\u00f73=98, 1\", \"332\u00f79=36, 8\", \"381\u00f75=76, 1\", \"606\u00f77=86, 4\", \"660\u00f76=110, 0\"],\n  8: [\"885\u00f74=221, 1\", \"594\u00f72=297, 0\", \"766\u00f76=127, 4\", \"949\u00f75=189, 4\", \"271\u00f77=38, 5\"],\n  12: [\"729\u00f76=121, 3\", \"710\u00f73=236, 2\", \"227\u00f76=37, 5\", \"780\u00f72=390, 0\", \"773\u00f78=96, 5\"],\n  16: [\"947\u00f73=315, 2\", \"820\u00f79=91, 1\", \"215\u00f72=107, 1\", \"313\u00f76=52, 1\", \"918\u00f79=102, 0\"],\n};\n\nfor (const rowIndexStr of Object.keys(updates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const values = updates[rowIndexStr];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const rng = cell.body.getRange();\n    rng.insertText(values[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet table of three-digit / one-digit division problems.\n# The table cell text is addressed positionally (row, column) rather than\n# via Find/Replace so that cells whose new value happens to equal another\n# cell's old value are not double-replaced.\n\n$d = $word.ActiveDocument\n$table = $d.Tables(1)\n\n# Map of 1-based row number -> array of new cell values (columns 1..5).\n$updates = @{\n    1  = @('382\u00f72=191, 0', '529\u00f78=66, 1', '704\u00f72=352, 0', '424\u00f78=53, 0', '145\u00f76=24, 1')\n    5  = @('295\u00f73=98, 1', '332\u00f79=36, 8', '381\u00f75=76, 1', '606\u00f77=86, 4', '660\u00f76=110, 0')\n    9  = @('885\u00f74=221, 1', '594\u00f72=297, 0', '766\u00f76=127, 4', '949\u00f75=189, 4', '271\u00f77=38, 5')\n    13 = @('729\u00f76=121, 3', '710\u00f73=236, 2', '227\u00f76=37, 5', '780\u00f72=390, 0', '773\u00f78=96, 5')\n    17 = @('947\u00f73=315, 2', '820\u00f79=91, 1', '215\u00f72=107, 1', '313\u00f76=52, 1', '918\u00f79=102, 0')\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $values = $updates[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $table.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
